$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header columns: "_old" -> "_FV2310", "_new" -> "_FV2404" ---
for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value2
    if ($val -ne $null) {
        if ($val.EndsWith("_old")) {
            $cell.Value = $val.Substring(0, $val.Length - 4) + "_FV2310"
        } elseif ($val.EndsWith("_new")) {
            $cell.Value = $val.Substring(0, $val.Length - 4) + "_FV2404"
        }
    }
}

# --- 2. Turn the used range into an Excel Table, preserving the header's
#        existing cell formatting without letting the table machinery bake
#        it into a dedicated header dxf (snapshot format to a scratch row,
#        strip it off the header before Add(), then restore afterwards). ---
$hdr = $ws.Range("A1:U1")
$scratch = $ws.Range("A1000:U1000")

$hdr.Copy()
$scratch.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$hdr.ClearFormats()

$rng = $ws.Range("A1:U79")
$listObj = $ws.ListObjects.Add(1, $rng, $false, 1, "")
$listObj.Name = "Table1"
$listObj.TableStyle = ""

$scratch.Copy()
$hdr.PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Rows.Item(1000).Delete()

# --- 3. Freeze the header row (pane split below row 1) ---
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
